$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Defs" (sheet1): refresh uid/created/updated for rows 2-4
# ---------------------------------------------------------------------------
$wsDefs = $wb.Worksheets.Item("Defs")

$wsDefs.Range("A2").Value = "lgs8j8fo-06px"
$wsDefs.Range("B2").Value = "2023-04-22T17:09:10.644Z"
$wsDefs.Range("C2").Value = "lgs8j8fo"

$wsDefs.Range("A3").Value = "lgs8j8fo-0oj7"
$wsDefs.Range("B3").Value = "2023-04-22T17:09:10.644Z"
$wsDefs.Range("C3").Value = "lgs8j8fo"

$wsDefs.Range("A4").Value = "lgs8j8fo-03ay"
$wsDefs.Range("B4").Value = "2023-04-22T17:09:10.644Z"
$wsDefs.Range("C4").Value = "lgs8j8fo"

# ---------------------------------------------------------------------------
# Sheet "Point Defs" (sheet2): insert a new "Select Test" row at row 2,
# refresh uid/created/updated on the (now shifted) rows, tweak a couple of
# values that changed for the existing rows
# ---------------------------------------------------------------------------
$wsPointDefs = $wb.Worksheets.Item("Point Defs")

$wsPointDefs.Rows.Item(2).Insert()

$wsPointDefs.Range("A2").Value = "lgs8j8fo-pjps"
$wsPointDefs.Range("B2").Value = "2023-04-22T17:09:10.644Z"
$wsPointDefs.Range("C2").Value = "lgs8j8fo"
$wsPointDefs.Range("D2").Value = $false
$wsPointDefs.Range("E2").Value = "0m7w"
$wsPointDefs.Range("F2").Value = "8esq"
$wsPointDefs.Range("G2").Value = "Select Test"
$wsPointDefs.Range("H2").Value = "⛏️"
$wsPointDefs.Range("I2").Value = "For testing selects"
$wsPointDefs.Range("J2").Value = "SELECT"
$wsPointDefs.Range("K2").Value = "COUNTOFEACH"

$wsPointDefs.Range("A3").Value = "lgs8j8fo-27z6i"
$wsPointDefs.Range("B3").Value = "2023-04-22T17:09:10.644Z"
$wsPointDefs.Range("C3").Value = "lgs8j8fo"
$wsPointDefs.Range("F3").Value = "1vb5"

$wsPointDefs.Range("A4").Value = "lgs8j8fo-62i6"
$wsPointDefs.Range("B4").Value = "2023-04-22T17:09:10.644Z"
$wsPointDefs.Range("C4").Value = "lgs8j8fo"
$wsPointDefs.Range("K4").Value = "AVERAGE"

$wsPointDefs.Range("A5").Value = "lgs8j8fo-nljl"
$wsPointDefs.Range("B5").Value = "2023-04-22T17:09:10.644Z"
$wsPointDefs.Range("C5").Value = "lgs8j8fo"

# ---------------------------------------------------------------------------
# Sheet "Entry Base" (sheet3): refresh uid/created/updated for rows 2-3,
# and the eid/period reference on row 3
# ---------------------------------------------------------------------------
$wsEntryBase = $wb.Worksheets.Item("Entry Base")

$wsEntryBase.Range("A2").Value = "lgs8j8fo-0lfw"
$wsEntryBase.Range("B2").Value = "2023-04-22T17:09:10.644Z"
$wsEntryBase.Range("C2").Value = "lgs8j8fo"

$wsEntryBase.Range("A3").Value = "lgs8j8fo-s0ps"
$wsEntryBase.Range("B3").Value = "2023-04-22T17:09:10.644Z"
$wsEntryBase.Range("C3").Value = "lgs8j8fo"
$wsEntryBase.Range("F3").Value = "lgs8j8g0-mpib"
$wsEntryBase.Range("G3").Value = "2023-04-22T12:09:10"

# ---------------------------------------------------------------------------
# Sheet "Entry Points" (sheet4): refresh uid/created/updated for rows 2-3
# ---------------------------------------------------------------------------
$wsEntryPoints = $wb.Worksheets.Item("Entry Points")

$wsEntryPoints.Range("A2").Value = "lgs8j8fo-afsz"
$wsEntryPoints.Range("B2").Value = "2023-04-22T17:09:10.644Z"
$wsEntryPoints.Range("C2").Value = "lgs8j8fo"

$wsEntryPoints.Range("A3").Value = "lgs8j8fo-x1oi"
$wsEntryPoints.Range("B3").Value = "2023-04-22T17:09:10.644Z"
$wsEntryPoints.Range("C3").Value = "lgs8j8fo"

# ---------------------------------------------------------------------------
# Sheet "Tag Defs" (sheet5): drop the _emoji/_desc header columns, add three
# new tag rows
# ---------------------------------------------------------------------------
$wsTagDefs = $wb.Worksheets.Item("Tag Defs")

$wsTagDefs.Range("G1:H1").ClearContents()

$wsTagDefs.Range("A2").Value = "lgs8j8g0-063q"
$wsTagDefs.Range("B2").Value = "2023-04-22T17:09:10.656Z"
$wsTagDefs.Range("C2").Value = "lgs8j8g0"
$wsTagDefs.Range("D2").Value = $false
$wsTagDefs.Range("E2").Value = "0q9d"
$wsTagDefs.Range("F2").Value = "My Tag!"

$wsTagDefs.Range("A3").Value = "lgs8j8g0-r9pi"
$wsTagDefs.Range("B3").Value = "2023-04-22T17:09:10.656Z"
$wsTagDefs.Range("C3").Value = "lgs8j8g0"
$wsTagDefs.Range("D3").Value = $false
$wsTagDefs.Range("E3").Value = "vvct"
$wsTagDefs.Range("F3").Value = "Orig Tag Label"

$wsTagDefs.Range("A4").Value = "lgs8j8g0-83ol"
$wsTagDefs.Range("B4").Value = "2023-04-22T17:09:10.656Z"
$wsTagDefs.Range("C4").Value = "lgs8j8g0"
$wsTagDefs.Range("D4").Value = $false
$wsTagDefs.Range("E4").Value = "0vvi"
$wsTagDefs.Range("F4").Value = "Select Option Test"

# ---------------------------------------------------------------------------
# Sheet "Tags" (sheet6): add two new tag-assignment rows
# ---------------------------------------------------------------------------
$wsTags = $wb.Worksheets.Item("Tags")

$wsTags.Range("A2").Value = "lgs8j8g0-31g6"
$wsTags.Range("B2").Value = "2023-04-22T17:09:10.656Z"
$wsTags.Range("C2").Value = "lgs8j8g0"
$wsTags.Range("D2").Value = $false
$wsTags.Range("E2").Value = "ay7l"
$wsTags.Range("G2").Value = "vvct"

$wsTags.Range("A3").Value = "lgs8j8g0-uq0p"
$wsTags.Range("B3").Value = "2023-04-22T17:09:10.656Z"
$wsTags.Range("C3").Value = "lgs8j8g0"
$wsTags.Range("D3").Value = $false
$wsTags.Range("E3").Value = "0m7w"
$wsTags.Range("F3").Value = "8esq"
$wsTags.Range("G3").Value = "0vvi"
